$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the price column headers to include the currency unit "(đồng)"
$ws.Range("G1").Value = "Giá nhập (đồng)"
$ws.Range("H1").Value = "Giá bán (đồng)"

# Move the active selection to H1 as in the edited workbook
$ws.Range("H1").Select()
